# Auto-generated Excel COM-interop edit script
# Applies scheduled-runner market-price refresh values to the Yojimbo_Profits workbook sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3519.186
$ws.Range("I40").Value = 4272.5938
$ws.Range("J40").Value = 1327.4546
$ws.Range("K40").Value = 4272.5938
$ws.Range("L40").Value = 1327.4546
$ws.Range("M40").Value = -4097.5938
$ws.Range("N40").Value = -1677.4546

$ws.Range("H69").Value = 4700
$ws.Range("J69").Value = 5000
$ws.Range("L69").Value = 15000
$ws.Range("N69").Value = -16748

$ws.Range("H72").Value = 4700
$ws.Range("J72").Value = 5000
$ws.Range("L72").Value = 45000
$ws.Range("N72").Value = -53736

$ws.Range("H76").Value = 3543.9443
$ws.Range("I76").Value = 3543.9443
$ws.Range("K76").Value = 3543.9443
$ws.Range("M76").Value = -3228.9443

$ws.Range("H79").Value = 3543.9443
$ws.Range("I79").Value = 3543.9443
$ws.Range("K79").Value = 3543.9443
$ws.Range("M79").Value = -2451.9443

$ws.Range("H80").Value = 1051
$ws.Range("I80").Value = 575
$ws.Range("J80").Value = 2003
$ws.Range("K80").Value = 1725
$ws.Range("L80").Value = 6009
$ws.Range("M80").Value = -727
$ws.Range("N80").Value = -8005

$ws.Range("H83").Value = 1051
$ws.Range("I83").Value = 575
$ws.Range("J83").Value = 2003
$ws.Range("K83").Value = 5175
$ws.Range("L83").Value = 18027
$ws.Range("M83").Value = -183
$ws.Range("N83").Value = -28011

$ws.Range("H121").Value = 835.25
$ws.Range("J121").Value = 825.6667
$ws.Range("L121").Value = 2477.0001
$ws.Range("N121").Value = -5971.0001

$ws.Range("H129").Value = 746.2632
$ws.Range("J129").Value = 868.3333
$ws.Range("L129").Value = 2604.9999
$ws.Range("N129").Value = -12604.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 30582.818
$ws.Range("J92").Value = 30582.818
$ws.Range("L92").Value = 30582.818
$ws.Range("N92").Value = -35574.818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 854
$ws.Range("I22").Value = 937.4286
$ws.Range("J22").Value = 708
$ws.Range("K22").Value = 937.4286
$ws.Range("L22").Value = 708
$ws.Range("M22").Value = -587.4286
$ws.Range("N22").Value = -1408

$ws.Range("H31").Value = 3843.1562
$ws.Range("I31").Value = 3187.1
$ws.Range("J31").Value = 4936.5835
$ws.Range("K31").Value = 3187.1
$ws.Range("L31").Value = 4936.5835
$ws.Range("M31").Value = -2892.1
$ws.Range("N31").Value = -5526.5835

$ws.Range("H34").Value = 3843.1562
$ws.Range("I34").Value = 3187.1
$ws.Range("J34").Value = 4936.5835
$ws.Range("K34").Value = 3187.1
$ws.Range("L34").Value = 4936.5835
$ws.Range("M34").Value = -2985.1
$ws.Range("N34").Value = -5340.5835

$ws.Range("H122").Value = 2718.6667
$ws.Range("I122").Value = 3525.1428
$ws.Range("J122").Value = 2124.4211
$ws.Range("K122").Value = 10575.4284
$ws.Range("L122").Value = 6373.263300000001
$ws.Range("M122").Value = -8125.428400000001
$ws.Range("N122").Value = -11273.2633

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 568.25
$ws.Range("I5").Value = 229.625
$ws.Range("J5").Value = 2600
$ws.Range("K5").Value = 688.875
$ws.Range("L5").Value = 7800
$ws.Range("M5").Value = -576.875
$ws.Range("N5").Value = -8024

$ws.Range("H109").Value = 3587.1853
$ws.Range("I109").Value = 972.4
$ws.Range("J109").Value = 5125.294
$ws.Range("K109").Value = 2917.2
$ws.Range("L109").Value = 15375.882
$ws.Range("M109").Value = -1877.2
$ws.Range("N109").Value = -17455.882

$ws.Range("H131").Value = 906.85
$ws.Range("J131").Value = 940.93475
$ws.Range("L131").Value = 2822.80425
$ws.Range("N131").Value = -12902.80425

$ws.Range("H135").Value = 568.25
$ws.Range("I135").Value = 229.625
$ws.Range("J135").Value = 2600
$ws.Range("K135").Value = 2066.625
$ws.Range("L135").Value = 23400
$ws.Range("M135").Value = 468.375
$ws.Range("N135").Value = -28470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1500
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -705

$ws.Range("H27").Value = 1500
$ws.Range("I27").Value = 1000
$ws.Range("K27").Value = 1000
$ws.Range("M27").Value = -893

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents() | Out-Null

$ws.Range("H132").Value = 6817.086
$ws.Range("I132").Value = 4220.625
$ws.Range("J132").Value = 12482.091
$ws.Range("K132").Value = 12661.875
$ws.Range("L132").Value = 37446.273
$ws.Range("M132").Value = -10131.875
$ws.Range("N132").Value = -42506.273

$ws.Range("H136").Value = 2747.4468
$ws.Range("I136").Value = 2073.4482
$ws.Range("J136").Value = 3833.3333
$ws.Range("K136").Value = 6220.344599999999
$ws.Range("L136").Value = 11499.9999
$ws.Range("M136").Value = -3670.344599999999
$ws.Range("N136").Value = -16599.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5343
$ws.Range("I81").Value = 1350.25
$ws.Range("J81").Value = 10666.667
$ws.Range("K81").Value = 2700.5
$ws.Range("L81").Value = 21333.334
$ws.Range("M81").Value = -1639.5
$ws.Range("N81").Value = -23455.334

$ws.Range("H84").Value = 5343
$ws.Range("I84").Value = 1350.25
$ws.Range("J84").Value = 10666.667
$ws.Range("K84").Value = 13502.5
$ws.Range("L84").Value = 106666.67
$ws.Range("M84").Value = -8198.5
$ws.Range("N84").Value = -117274.67

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents() | Out-Null

$ws.Range("H132").Value = 1138.0294
$ws.Range("I132").Value = 786.3043
$ws.Range("J132").Value = 1873.4546
$ws.Range("K132").Value = 2358.9129
$ws.Range("L132").Value = 5620.3638
$ws.Range("M132").Value = 171.0870999999997
$ws.Range("N132").Value = -10680.3638

$ws.Range("H136").Value = 756.3889
$ws.Range("I136").Value = 652.4194
$ws.Range("J136").Value = 1401
$ws.Range("K136").Value = 1957.2582
$ws.Range("L136").Value = 4203
$ws.Range("M136").Value = 592.7418
$ws.Range("N136").Value = -9303
